# 4-state model for spring females: add the 4-state transition header row
# (row 17) and an "NA"/-1000 placeholder data row (row 18), reusing the
# existing "n -> m" / "NA" shared strings already used by the 5-state
# block in rows 14-15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 - transition-label header row (only the 4-state columns are filled;
# columns M-P, which only make sense in the 5-state model, are left blank)
$ws.Range("A17").Value = "1 -> 2"
$ws.Range("B17").Value = "1 -> 3"
$ws.Range("C17").Value = "1 -> 4"
$ws.Range("D17").Value = "2 -> 1"
$ws.Range("E17").Value = "2 -> 3"
$ws.Range("F17").Value = "2 -> 4"
$ws.Range("G17").Value = "3 -> 1"
$ws.Range("H17").Value = "3 -> 2"
$ws.Range("I17").Value = "3 -> 4"
$ws.Range("J17").Value = "4 -> 1"
$ws.Range("K17").Value = "4 -> 2"
$ws.Range("L17").Value = "4 -> 3"
$ws.Range("Q17").Value = "5 -> 1"
$ws.Range("R17").Value = "5 -> 2"
$ws.Range("S17").Value = "5 -> 3"
$ws.Range("T17").Value = "5 -> 4"

# Row 18 - placeholder values: -1000 for numeric (estimable) transitions,
# "NA" for the transitions that are structurally impossible in this model
$ws.Range("A18").Value = "NA"
$ws.Range("B18").Value = -1000
$ws.Range("C18").Value = -1000
$ws.Range("D18").Value = -1000
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = -1000
$ws.Range("H18").Value = "NA"
$ws.Range("I18").Value = -1000
$ws.Range("J18").Value = -1000
$ws.Range("K18").Value = -1000
$ws.Range("L18").Value = -1000
$ws.Range("Q18").Value = -1000
$ws.Range("R18").Value = -1000
$ws.Range("S18").Value = -1000
$ws.Range("T18").Value = "NA"

# Scroll the view so column B is at the left edge, and select the new row
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B18:L18").Select()

$wb.Save()
